$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Validation")

$ws.Range("A3").Value = 40901
$ws.Range("B3").Value = 46627
$ws.Range("C3").Value = 51326
$ws.Range("D3").Value = 65202
$ws.Range("I3").Value = 1594
$ws.Range("J3").Value = 1775
$ws.Range("K3").Value = 1910
$ws.Range("L3").Value = 2347
